$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a value as literal TEXT (not auto-converted to a number),
# while leaving the destination cell's existing style/format untouched.
# We stage the text in a scratch cell that is forced to Text format ("@"),
# copy it, and paste-special VALUES ONLY into the destination so the
# destination keeps its original cell style index.
# ---------------------------------------------------------------------------
function Set-TextValue($rangeAddr, $text) {
    $scratch = $ws.Range("ZZ1000")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($rangeAddr).PasteSpecial(-4163) # xlPasteValues
    $scratch.Clear()
}

# ---------------------------------------------------------------------------
# Top summary metrics (K/L columns near the top of the sheet)
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 67
$ws.Range("L7").Value = 5
Set-TextValue "L9"  "26.0%"
Set-TextValue "L10" "77.4%"

# ---------------------------------------------------------------------------
# Group-statistics block (rows 15, 18-23): Missing (O) / recorded count (P)
# and the derived percentages (R / S) shift for every B1-x group except B1-11
# (row 17, handled separately below because its whole row also changes
# recorded-status).
# ---------------------------------------------------------------------------
$groupRows = @(15,18,19,20,21,22,23)
foreach ($r in $groupRows) {
    $ws.Range("O$r").Value = 6
    $ws.Range("P$r").Value = 0
}

Set-TextValue "R15" "28.6%"
Set-TextValue "S15" "84.6%"

Set-TextValue "R18" "27.3%"
Set-TextValue "S18" "74.6%"

Set-TextValue "R19" "28.6%"
Set-TextValue "S19" "81.2%"

Set-TextValue "R20" "28.6%"
Set-TextValue "S20" "79.2%"

Set-TextValue "R21" "28.6%"
Set-TextValue "S21" "78.7%"

Set-TextValue "R22" "28.6%"
Set-TextValue "S22" "77.3%"

Set-TextValue "R23" "28.6%"
Set-TextValue "S23" "73.3%"

# S17 (average attendance % for the B1-11 group) also changes even though
# the row itself keeps its own group stats besides this cell.
Set-TextValue "S17" "71.6%"

# ---------------------------------------------------------------------------
# Session rows that flip from "Not Recorded" to "Recorded" because the
# attendance system has now picked up a count of students for them. These
# rows copy the formatting already used by other "Recorded" rows (e.g. A23)
# so the green highlight/font match exactly, then get their
# Recorded-By / Students / Status cells updated.
# ---------------------------------------------------------------------------
function Set-RecordedRow($rowNum, $studentsText) {
    $srcFormat = $ws.Range("A23:I23")
    $dstFormat = $ws.Range("A" + $rowNum + ":I" + $rowNum)
    $srcFormat.Copy()
    $dstFormat.PasteSpecial(-4122) # xlPasteFormats
    $ws.Range("G$rowNum").Value = "System"
    Set-TextValue "H$rowNum" $studentsText
    $ws.Range("I$rowNum").Value = "Recorded"
}

Set-RecordedRow 17  "21/27"
Set-RecordedRow 83  "2/21"
Set-RecordedRow 104 "1/31"
Set-RecordedRow 125 "1/28"
Set-RecordedRow 146 "1/29"
Set-RecordedRow 167 "2/33"
Set-RecordedRow 188 "1/30"

# ---------------------------------------------------------------------------
# Straightforward student-count text update (no status/format change)
# ---------------------------------------------------------------------------
Set-TextValue "H46" "8/19"

$excel.CutCopyMode = 0
